$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (I1, J1) - match the styling of the existing header row by copying H1's style
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$data = @(
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(10, 10),
    @(8, 8),
    @(4, 6),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(5, 7),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(5, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
